# Update countries & provincias Spain
# Applies the 14-Abril-2020 17:22 COVID-19 data refresh to the "Pais" sheet:
#  - bumps the "Datos actualizados" timestamp in A1
#  - updates case counters for countries whose totals changed
#  - re-sorts (by swapping row contents) wherever the new totals change the
#    descending ranking, exactly as the source diff shows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 17:22"

# --- helper: write one full data row (A..H) ---------------------------
function Set-PaisRow {
    param($row, $pais, $totales, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes)

    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $totales
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Alemania: stats refresh only, same rank (row 8)
Set-PaisRow 8 "Alemania" 130694 622 68200 59233 4288 67 3261

# Chile overtakes Japon and Ecuador (rows 27-29 shift down)
Set-PaisRow 27 "Chile"   7917 392 2646 5179 387 10 92
Set-PaisRow 28 "Japon"   7645 27  799  6703 135 0  143
Set-PaisRow 29 "Ecuador" 7529 0   597  6577 121 0  355

# Grecia: stats refresh only, same rank (row 57)
Set-PaisRow 57 "Grecia" 2170 25 269 1800 73 2 101

# Cuba overtakes Tunez (rows 82-83 shift down)
Set-PaisRow 82 "Cuba"  766 40 132 613 11 0 21
Set-PaisRow 83 "Tunez" 726 0  43  649 89 0 34

# Republica de Chipre: stats refresh only, same rank (row 87)
Set-PaisRow 87 "Republica de Chipre" 695 33 65 618 8 0 12

# Mauricio: stats refresh only, same rank (row 109)
Set-PaisRow 109 "Mauricio" 324 0 51 264 3 0 9

# Birmania: stats refresh only, same rank (row 142)
Set-PaisRow 142 "Birmania" 63 1 2 57 0 0 4

# Benin: stats refresh only, same rank (row 160)
Set-PaisRow 160 "Benin" 35 0 18 16 0 0 1

# Siria jumps ahead of Sudan, Mozambique, Libia (rows 165-170 shift down)
Set-PaisRow 165 "Siria"               29 4 5 22 0 0 2
Set-PaisRow 166 "Sudan"               29 0 4 21 0 0 4
Set-PaisRow 167 "Mozambique"          28 7 2 26 0 0 0
Set-PaisRow 168 "Libia"               26 0 9 16 0 0 1
Set-PaisRow 169 "Republica del Chad"  23 0 2 21 0 0 0
Set-PaisRow 170 "Antigua y Barbuda"   23 0 3 18 1 0 2
